# Issue 50917: imported folder archives do not correct file paths for file
# fields - selenium tests.
#
# The GenericAssay_Run2.xlsx run-data sheet gets a new "resultFileField"
# column (F) with a sample file-path value ("help.jpg") used by the assay
# import/export selenium tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column F.
$ws.Range("F1").Value = "resultFileField"

# Sample file reference for the first data row.
$ws.Range("F2").Value = "help.jpg"

# Leave the cursor where the user would after typing the value into F2.
$ws.Range("F3").Select() | Out-Null
